$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 71 ("co2_per_egg" entry) pushing the crop
# rows (previously 71-86) down to 72-87, and fill it with the new data
# belonging to the CO2 emissions group (matching the style used by the
# existing CO2 rows 67-70).
$ws.Rows.Item(71).Insert()

# Match the formatting/style of the existing CO2 rows (row 70) so the new
# row belongs visually to the same group.
$ws.Range("A70:F70").Copy()
$ws.Range("A71:F71").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A71").Value = "co2_per_egg"
$ws.Range("B71").Value = 0.1
$ws.Range("C71").Value = "NA"
$ws.Range("D71").Value = 0.16
$ws.Range("E71").Value = "posnorm"
$ws.Range("F71").Value = "CO2 per Egg"

$ws.Range("A71:F71").Select()

$ws.Range("F79").Select()
